$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 129; this shifts the existing rows 129:187
# down to 130:188 (carrying all their data/styles with them).
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new record.
$ws.Range("A129").Value = 8
$ws.Range("B129").Value = "Terminal La Palmera de La Serena"
$ws.Range("C129").Value = "Coquimbo"
$ws.Range("D129").Value = 44596
$ws.Range("E129").Value = 4
$ws.Range("F129").Value = 100112021
$ws.Range("G129").Value = "Ají"
$ws.Range("H129").Value = "Americana (o)"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 560
$ws.Range("K129").Value = 11500
$ws.Range("L129").Value = 12000
$ws.Range("M129").Value = 11750
$ws.Range("N129").Value = "$/caja 15 kilos"
$ws.Range("O129").Value = "Provincia de Limarí"
$ws.Range("P129").Value = 783
$ws.Range("Q129").Value = 15
$ws.Range("R129").Value = "Hortaliza"
